# "fixed the scale for hair" - update the Mark Hair (column C) scores on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of cell -> corrected value (Mark Hair scale fix)
$updates = @{
    "C13"  = 0
    "C16"  = 0
    "C21"  = 2
    "C26"  = 2
    "C28"  = 2
    "C29"  = 3
    "C35"  = 1
    "C38"  = 2
    "C43"  = 3
    "C48"  = 3
    "C50"  = 2
    "C51"  = 2
    "C54"  = 2
    "C56"  = 3
    "C76"  = 2
    "C77"  = 3
    "C90"  = 2
    "C92"  = 2
    "C94"  = 3
    "C95"  = 3
    "C96"  = 2
    "C100" = 3
    "C101" = 2
    "C113" = 3
    "C115" = 2
    "C117" = 3
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# Leave the sheet scrolled/selected where the author ended up working.
$ws.Activate()
$ws.Range("F110").Select()
